$p = $ppt.ActivePresentation

# --- Update the watermark text on slide 1 (TextBox shape): bump GroupDocs.Assembly version ---
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Remember the shape's current geometry (points) so it can be restored after the text
# edit - the shape auto-fits its text, and touching the run text recomputes the
# bounding box; the version bump itself is not supposed to resize/move the shape.
$origLeft = [double]$shp.Left
$origTop = [double]$shp.Top
$origWidth = [double]$shp.Width
$origHeight = [double]$shp.Height

$full = $tr.Text
$oldText = "Created with GroupDocs.Assembly 25.6."
$newText = "Created with GroupDocs.Assembly 25.12."
$startIdx = $full.IndexOf($oldText)
if ($startIdx -ge 0) {
    $sub = $tr.Characters($startIdx + 1, $oldText.Length)
    $sub.Text = $newText
}

# Restore the original geometry. A tiny relative nudge compensates for float32
# round-tripping through the Left/Top/Width/Height properties so the restored
# values land back on the exact original EMU values.
$eps = 0.00000015
$shp.Left = $origLeft * (1 + $eps)
$shp.Top = $origTop * (1 + $eps)
$shp.Width = $origWidth * (1 + $eps)
$shp.Height = $origHeight * (1 + $eps)

# --- Update presentation tags (AS_NET, AS_OS, AS_TITLE) ---
$p.Tags.Add("AS_NET", "4.0.30319.42000")
$p.Tags.Add("AS_OS", "Microsoft Windows NT 6.2.9200.0")
$p.Tags.Add("AS_TITLE", "Aspose.Slides for .NET 4.0 Client Profile")
